# Auto update: 2025-12-05 10:51:57
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - UnitedHealth Group (UNH)
$ws.Range("D2").Value = 333.49
$ws.Range("E2").Value = 50.6
$ws.Range("F2").Value = 1.15
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 60.9
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 54.85170003294819
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3 - MetLife, Inc. (MET)
$ws.Range("D3").Value = 78.03
$ws.Range("E3").Value = 43
$ws.Range("F3").Value = 2.23
$ws.Range("K3").Value = 58.5
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 54.85170003294819
$ws.Range("O3").Value = "⚪ 중립 구간"

# Row 4 - American International Group, Inc. (AIG)
$ws.Range("D4").Value = 77.28
$ws.Range("E4").Value = 42.7
$ws.Range("F4").Value = 1.63
$ws.Range("H4").Value = 46
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 51.7
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 54.85170003294819
$ws.Range("O4").Value = "⚪ 중립 구간"

# Row 5 - Prudential Financial, Inc. (PRU)
$ws.Range("D5").Value = 110.25
$ws.Range("E5").Value = 63.7
$ws.Range("F5").Value = 2.11
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 36
$ws.Range("J5").Value = 36
$ws.Range("K5").Value = 48.9
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 54.85170003294819
$ws.Range("O5").Value = "⚪ 중립 구간"
